$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.274.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5235"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06325"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.665.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.893.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5471"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8282"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.301.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.677"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.069"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.197"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06173"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.281"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.591"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.294"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9714"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.427"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5751"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01609"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.006"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.023.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.809.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈110"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.063"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.487"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "
